$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44300
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 2143

# Row 3
$ws.Range("D3").Value = 44300
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("S3").Value = 1714

# Row 4
$ws.Range("D4").Value = 44320
$ws.Range("M4").Value = 20
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 1714

# Row 5
$ws.Range("D5").Value = 44320
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 1143

# Row 6
$ws.Range("D6").Value = 44292
$ws.Range("M6").Value = 25
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("S6").Value = 2286

# Row 7
$ws.Range("D7").Value = 44292
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("S7").Value = 2143

# Row 8
$ws.Range("D8").Value = 44301
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 14000
$ws.Range("P8").Value = 14000
$ws.Range("S8").Value = 2000

# Row 9
$ws.Range("D9").Value = 44301
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("S9").Value = 1714

# Row 10
$ws.Range("D10").Value = 44322
$ws.Range("M10").Value = 45
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 12000
$ws.Range("S10").Value = 1714

# Row 11
$ws.Range("D11").Value = 44322
$ws.Range("N11").Value = 8000
$ws.Range("O11").Value = 8000
$ws.Range("P11").Value = 8000
$ws.Range("S11").Value = 1143

# Row 12
$ws.Range("D12").Value = 44299
$ws.Range("M12").Value = 80
$ws.Range("R12").Value = "Provincia de Santiago"

# Row 13
$ws.Range("D13").Value = 44299
$ws.Range("M13").Value = 75
$ws.Range("R13").Value = "Provincia de Santiago"

# Row 14
$ws.Range("D14").Value = 44302
$ws.Range("M14").Value = 50

# Row 15
$ws.Range("D15").Value = 44302
$ws.Range("M15").Value = 30
